$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 4, 5) are being cyclically rotated:
#   new Row2 <- old Row5
#   new Row4 <- old Row2
#   new Row5 <- old Row4
# Row 3 is untouched.
# Capture the "old" values for columns D and J..Q before overwriting anything.

$cols = @("D","J","K","L","M","N","O","P","Q")

$old2 = @{}
$old4 = @{}
$old5 = @{}

foreach ($col in $cols) {
    $old2[$col] = $ws.Range("${col}2").Value()
    $old4[$col] = $ws.Range("${col}4").Value()
    $old5[$col] = $ws.Range("${col}5").Value()
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $old5[$col]
    $ws.Range("${col}4").Value = $old2[$col]
    $ws.Range("${col}5").Value = $old4[$col]
}
